$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Main input: lower the pizza price, which ripples through every
# formula-driven cell (row 9-14 incl. the revenue table) and the
# two charts that read from them.
$ws.Range("F6").Value = 4

# New "optimal mix" labels next to the revenue-maximizing OFFSET lookups.
$ws.Range("A13").Copy()
$ws.Range("M12").PasteSpecial(-4122)
$ws.Range("M12").Value = "Opt. pizza:"

$ws.Range("A13").Copy()
$ws.Range("M13").PasteSpecial(-4122)
$ws.Range("M13").Value = "Opt. book:"

# Make the highlighted (revenue-maximizing) point on the PPF chart bigger.
$chart = $ws.ChartObjects(1).Chart
$highlight = $chart.SeriesCollection("Highlight")
$highlight.MarkerSize = 6

$ws.Range("F7").Select()
